# Commit: "update file system typos"
# Adds a new observation row (row 5) with a value of 8 in column C,
# and leaves the selection on D6 (matching the saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = 8

$ws.Range("D6").Select()
